# BRUTO Y TARA modification / CUPOS new version
# Populates the "eje" (axle) reference table in columns D:F, rows 42-50,
# and fills the O column (COD RUB counts) for rows 48-62 in the
# RUBROSXPRODUCTO.DAT table. Also updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New small lookup table in D:F (rows 42-50) -----------------------
# Order below intentionally mirrors how the values were entered so the
# resulting shared-string table matches the source order.
$ws.Range("F42").Value = "50-60"
$ws.Range("E42").Value = "ilegal"

$ws.Range("E43").Value = "lega"
$ws.Range("F43").Value = 45

$ws.Range("E44").Value = "tara"
$ws.Range("F44").Value = 30

$ws.Range("E46").Value = "ejes"

$ws.Range("E47").Value = "simple"
$ws.Range("F47").Value = 10.5

$ws.Range("E48").Value = "doble "
$ws.Range("D48").Value = "tandem "
$ws.Range("F48").Value = 18

$ws.Range("E49").Value = "triple"
$ws.Range("F49").Value = 25

$ws.Range("D50").Value = "rueda "
$ws.Range("E50").Value = "direccional"
$ws.Range("F50").Value = 6

$ws.Range("D45").Value = "peso cami"
$ws.Range("E45").Value = 8
$ws.Range("F45").Value = 15

# --- New counts in column O for the RUBROSXPRODUCTO.DAT block --------
$ws.Range("O48").Value = 1
$ws.Range("O49").Value = 2
$ws.Range("O50").Value = 4
$ws.Range("O51").Value = 1
$ws.Range("O52").Value = 4
$ws.Range("O53").Value = 5
$ws.Range("O54").Value = 1
$ws.Range("O55").Value = 5
$ws.Range("O56").Value = 9
$ws.Range("O57").Value = 1
$ws.Range("O58").Value = 4
$ws.Range("O59").Value = 8
$ws.Range("O60").Value = 1
$ws.Range("O61").Value = 4
$ws.Range("O62").Value = 5

# --- Update the active selection shown when the workbook is reopened -
$ws.Range("F44").Select()
